# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns
# for the crypto ranking sheet, matching the latest scrape.
# Price column cells are kept as plain text ("@" number format) since the
# source values (e.g. "2.227.19") are display strings, not real numbers,
# and some updated values (e.g. "242.89") would otherwise be
# auto-recognized by Excel as numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextCell "D2" "42.167.65"
Set-TextCell "E2" "  -0.77%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.227.19"
Set-TextCell "E3" "  -0.77%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.20%  "

# Row 5 - BNB
Set-TextCell "D5" "242.89"
Set-TextCell "E5" "  -0.92%  "

# Row 6 - XRP
Set-TextCell "E6" "  +0.88%  "

# Row 7 - Solana
Set-TextCell "D7" "74.23"
Set-TextCell "E7" "  -1.97%  "

# Row 8 - USDC
Set-TextCell "E8" "  +0.18%  "

# Row 9 - Cardano
Set-TextCell "E9" "  -3.10%  "

# Row 10 - Avalanche
Set-TextCell "D10" "42.72"
Set-TextCell "E10" "  -2.98%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0960"
Set-TextCell "E11" "  +1.19%  "

# Row 12 - Polkadot
Set-TextCell "E12" "  -3.43%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell "D14" "2.560.44"
Set-TextCell "E14" "  -0.63%  "

# Row 15 - Chainlink
Set-TextCell "E15" "  -1.53%  "

# Row 16 - Polygon
Set-TextCell "D16" "0.836"
Set-TextCell "E16" "  -2.83%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.232.41"
Set-TextCell "E17" "  -1.58%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "41.997.33"
Set-TextCell "E18" "  -0.69%  "

# Row 19 - ShibaInu
Set-TextCell "E19" "  +3.70%  "

# Row 20 - Uniswap
Set-TextCell "D20" "6.22"
Set-TextCell "E20" "  +0.15%  "

# Row 21 - Litecoin
Set-TextCell "D21" "72.83"
Set-TextCell "E21" "  +1.11%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextCell "D22" "11.04"
Set-TextCell "E22" "  -1.16%  "

# Row 23 - BitcoinCash
Set-TextCell "D23" "230.51"
Set-TextCell "E23" "  -0.55%  "

# Row 24 - ImmutableX
Set-TextCell "D24" "2.10"
Set-TextCell "E24" "  -5.85%  "

# Row 25 - Dai
Set-TextCell "E25" "  +0.23%  "

# Row 26 - Cosmos
Set-TextCell "E26" "  -2.95%  "

# Row 27 - WEMIXToken
Set-TextCell "E27" "  -0.11%  "

# Row 28 - PancakeSwap
Set-TextCell "D28" "2.28"
Set-TextCell "E28" "  -0.96%  "

# Row 29 - Toncoin
Set-TextCell "E29" "  -2.73%  "

# Row 30 - Monero
Set-TextCell "D30" "167.02"
Set-TextCell "E30" "  -0.01%  "

# Row 31 - EthereumClassic
Set-TextCell "D31" "20.61"
Set-TextCell "E31" "  -0.49%  "

# Row 32 - Filecoin
Set-TextCell "E32" "  -5.93%  "

# Row 33 - Hedera
Set-TextCell "D33" "0.0804"
Set-TextCell "E33" "  -1.52%  "

# Row 34 - InjectiveProtocol
Set-TextCell "D34" "30.02"
Set-TextCell "E34" "  -2.51%  "

# Row 35 - Stellar
Set-TextCell "E35" "  -0.55%  "

# Row 36 - Kaspa
Set-TextCell "E36" "  -7.62%  "

# Row 37 - RenderToken
Set-TextCell "D37" "4.33"
Set-TextCell "E37" "  -6.52%  "

# Row 38 - VeChain
Set-TextCell "E38" "  -3.68%  "

# Row 39 - Celestia
Set-TextCell "D39" "13.21"
Set-TextCell "E39" "  -4.24%  "

# Row 40 - LidoDAOToken
Set-TextCell "E40" "  -2.35%  "

# Row 41 - MultiversX
Set-TextCell "D41" "65.10"
Set-TextCell "E41" "  +2.03%  "

# Row 42 - THORChain
Set-TextCell "E42" "  -1.01%  "

# Row 43 - Algorand
Set-TextCell "E43" "  -0.40%  "

# Row 44 - FraxShare
Set-TextCell "D44" "8.72"
Set-TextCell "E44" "  -1.49%  "

# Row 45 - Aave
Set-TextCell "D45" "104.41"
Set-TextCell "E45" "  -2.33%  "

# Row 46 - Cronos
Set-TextCell "E46" "  -1.98%  "

# Row 47 - NEARProtocol
Set-TextCell "E47" "  -2.96%  "

# Row 48 - ARBITRUM
Set-TextCell "E48" "  -1.84%  "

# Row 49 - TrustWalletToken
Set-TextCell "E49" "  -1.08%  "

# Row 50 - HuobiToken
Set-TextCell "E50" "  -1.37%  "

# Row 51 - RocketPoolETH
Set-TextCell "D51" "2.430.91"
Set-TextCell "E51" "  -0.96%  "
